$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the duplicated "largest accepted serial number" row (old row 11,
#    value 2958465) -- its data already lives in row 10 after this delete,
#    shifting every row below it up by one.
# ---------------------------------------------------------------------------
$ws.Rows.Item(11).Delete()

# ---------------------------------------------------------------------------
# 2. Update the "every day" example block (now rows 7-10) with the new
#    serial-number values used in the refreshed example.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 46016
$ws.Range("A10").Value = 2958465

# Row 10 now carries the "largest accepted serial number" comment that used
# to live on the row we just removed (and drops the placeholder formatting
# that the empty D10 used to have).
$ws.Range("D10").ClearFormats()
$ws.Range("D10").Value = "Largest accepted serial number"

# ---------------------------------------------------------------------------
# 3. Insert a brand new row for the "zero serial number" example, right
#    before the "serial number too large" row (which is now row 18).
# ---------------------------------------------------------------------------
$ws.Rows.Item(18).Insert()

$ws.Range("A18").Formula = "=DAY(0)"
$ws.Range("B18").Formula = "=FORMULATEXT(A18)"
$ws.Range("C18").Value = "Zero serial number"

# ---------------------------------------------------------------------------
# 4. The former "negative serial number" example (now row 17) switches from
#    DAY(-1) to DAY(-5); the explanatory comment stays the same.
# ---------------------------------------------------------------------------
$ws.Range("A17").Formula = "=DAY(-5)"

# ---------------------------------------------------------------------------
# 5. Cosmetic: move the active selection to B12, matching where the author
#    last clicked while editing.
# ---------------------------------------------------------------------------
$ws.Range("B12").Select()
